$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (bold/border style for col A, date format for col B) to new rows 19:37
# by copying the format from the last existing formatted row (row 18).
$ws.Range("A18:H18").Copy()
$ws.Range("A19:H37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Build full data block for A2:H37 (36 rows x 8 cols)
$arr = New-Object "object[,]" 36,8
$arr[0,0] = 0
$arr[0,1] = 45271.30555555555
$arr[0,2] = 2
$arr[0,3] = 63
$arr[0,4] = 0
$arr[0,5] = 27
$arr[0,6] = 0
$arr[0,7] = 9
$arr[1,0] = 1
$arr[1,1] = 45271.3125
$arr[1,2] = 1
$arr[1,3] = 96
$arr[1,4] = 1
$arr[1,5] = 24
$arr[1,6] = 0
$arr[1,7] = 12
$arr[2,0] = 2
$arr[2,1] = 45271.31944444445
$arr[2,2] = 0
$arr[2,3] = 88
$arr[2,4] = 0
$arr[2,5] = 29
$arr[2,6] = 1
$arr[2,7] = 25
$arr[3,0] = 3
$arr[3,1] = 45271.32638888889
$arr[3,2] = 1
$arr[3,3] = 83
$arr[3,4] = 1
$arr[3,5] = 28
$arr[3,6] = 1
$arr[3,7] = 32
$arr[4,0] = 4
$arr[4,1] = 45271.33333333334
$arr[4,2] = 1
$arr[4,3] = 78
$arr[4,4] = 1
$arr[4,5] = 37
$arr[4,6] = 2
$arr[4,7] = 31
$arr[5,0] = 5
$arr[5,1] = 45271.34027777778
$arr[5,2] = 0
$arr[5,3] = 76
$arr[5,4] = 2
$arr[5,5] = 27
$arr[5,6] = 2
$arr[5,7] = 21
$arr[6,0] = 6
$arr[6,1] = 45271.34722222222
$arr[6,2] = 1
$arr[6,3] = 80
$arr[6,4] = 2
$arr[6,5] = 29
$arr[6,6] = 2
$arr[6,7] = 14
$arr[7,0] = 7
$arr[7,1] = 45271.35416666666
$arr[7,2] = 0
$arr[7,3] = 88
$arr[7,4] = 1
$arr[7,5] = 20
$arr[7,6] = 1
$arr[7,7] = 15
$arr[8,0] = 8
$arr[8,1] = 45271.36111111111
$arr[8,2] = 1
$arr[8,3] = 67
$arr[8,4] = 1
$arr[8,5] = 33
$arr[8,6] = 1
$arr[8,7] = 21
$arr[9,0] = 9
$arr[9,1] = 45271.36805555555
$arr[9,2] = 0
$arr[9,3] = 63
$arr[9,4] = 3
$arr[9,5] = 21
$arr[9,6] = 1
$arr[9,7] = 25
$arr[10,0] = 10
$arr[10,1] = 45271.375
$arr[10,2] = 0
$arr[10,3] = 83
$arr[10,4] = 3
$arr[10,5] = 22
$arr[10,6] = 2
$arr[10,7] = 9
$arr[11,0] = 11
$arr[11,1] = 45377.65972222222
$arr[11,2] = 3
$arr[11,3] = 77
$arr[11,4] = 0
$arr[11,5] = 23
$arr[11,6] = 0
$arr[11,7] = 14
$arr[12,0] = 12
$arr[12,1] = 45377.66666666666
$arr[12,2] = 2
$arr[12,3] = 63
$arr[12,4] = 0
$arr[12,5] = 25
$arr[12,6] = 2
$arr[12,7] = 16
$arr[13,0] = 13
$arr[13,1] = 45377.67361111111
$arr[13,2] = 2
$arr[13,3] = 48
$arr[13,4] = 0
$arr[13,5] = 16
$arr[13,6] = 1
$arr[13,7] = 23
$arr[14,0] = 14
$arr[14,1] = 45377.68055555555
$arr[14,2] = 1
$arr[14,3] = 60
$arr[14,4] = 0
$arr[14,5] = 20
$arr[14,6] = 1
$arr[14,7] = 22
$arr[15,0] = 15
$arr[15,1] = 45377.6875
$arr[15,2] = 2
$arr[15,3] = 67
$arr[15,4] = 0
$arr[15,5] = 26
$arr[15,6] = 1
$arr[15,7] = 24
$arr[16,0] = 16
$arr[16,1] = 45377.69444444445
$arr[16,2] = 3
$arr[16,3] = 82
$arr[16,4] = 0
$arr[16,5] = 13
$arr[16,6] = 2
$arr[16,7] = 24
$arr[17,0] = 17
$arr[17,1] = 45377.70138888889
$arr[17,2] = 2
$arr[17,3] = 63
$arr[17,4] = 0
$arr[17,5] = 19
$arr[17,6] = 3
$arr[17,7] = 26
$arr[18,0] = 18
$arr[18,1] = 45391.33333333334
$arr[18,2] = 3
$arr[18,3] = 78
$arr[18,4] = 0
$arr[18,5] = 41
$arr[18,6] = 0
$arr[18,7] = 25
$arr[19,0] = 19
$arr[19,1] = 45391.34027777778
$arr[19,2] = 0
$arr[19,3] = 55
$arr[19,4] = 1
$arr[19,5] = 31
$arr[19,6] = 2
$arr[19,7] = 36
$arr[20,0] = 20
$arr[20,1] = 45391.34722222222
$arr[20,2] = 3
$arr[20,3] = 76
$arr[20,4] = 0
$arr[20,5] = 21
$arr[20,6] = 1
$arr[20,7] = 35
$arr[21,0] = 21
$arr[21,1] = 45391.35416666666
$arr[21,2] = 1
$arr[21,3] = 78
$arr[21,4] = 0
$arr[21,5] = 31
$arr[21,6] = 1
$arr[21,7] = 21
$arr[22,0] = 22
$arr[22,1] = 45391.36111111111
$arr[22,2] = 1
$arr[22,3] = 53
$arr[22,4] = 3
$arr[22,5] = 30
$arr[22,6] = 2
$arr[22,7] = 27
$arr[23,0] = 23
$arr[23,1] = 45391.36805555555
$arr[23,2] = 4
$arr[23,3] = 59
$arr[23,4] = 1
$arr[23,5] = 26
$arr[23,6] = 3
$arr[23,7] = 24
$arr[24,0] = 24
$arr[24,1] = 45391.66666666666
$arr[24,2] = 1
$arr[24,3] = 103
$arr[24,4] = 2
$arr[24,5] = 21
$arr[24,6] = 1
$arr[24,7] = 16
$arr[25,0] = 25
$arr[25,1] = 45391.67361111111
$arr[25,2] = 2
$arr[25,3] = 84
$arr[25,4] = 1
$arr[25,5] = 19
$arr[25,6] = 1
$arr[25,7] = 16
$arr[26,0] = 26
$arr[26,1] = 45391.68055555555
$arr[26,2] = 2
$arr[26,3] = 70
$arr[26,4] = 3
$arr[26,5] = 21
$arr[26,6] = 1
$arr[26,7] = 19
$arr[27,0] = 27
$arr[27,1] = 45391.6875
$arr[27,2] = 1
$arr[27,3] = 98
$arr[27,4] = 1
$arr[27,5] = 44
$arr[27,6] = 2
$arr[27,7] = 13
$arr[28,0] = 28
$arr[28,1] = 45391.69444444445
$arr[28,2] = 1
$arr[28,3] = 73
$arr[28,4] = 3
$arr[28,5] = 17
$arr[28,6] = 0
$arr[28,7] = 16
$arr[29,0] = 29
$arr[29,1] = 45391.70138888889
$arr[29,2] = 4
$arr[29,3] = 88
$arr[29,4] = 0
$arr[29,5] = 19
$arr[29,6] = 3
$arr[29,7] = 12
$arr[30,0] = 30
$arr[30,1] = 45391.83333333334
$arr[30,2] = 1
$arr[30,3] = 69
$arr[30,4] = 0
$arr[30,5] = 14
$arr[30,6] = 2
$arr[30,7] = 9
$arr[31,0] = 31
$arr[31,1] = 45391.84027777778
$arr[31,2] = 0
$arr[31,3] = 35
$arr[31,4] = 0
$arr[31,5] = 10
$arr[31,6] = 1
$arr[31,7] = 11
$arr[32,0] = 32
$arr[32,1] = 45391.84722222222
$arr[32,2] = 1
$arr[32,3] = 46
$arr[32,4] = 2
$arr[32,5] = 5
$arr[32,6] = 0
$arr[32,7] = 9
$arr[33,0] = 33
$arr[33,1] = 45391.85416666666
$arr[33,2] = 0
$arr[33,3] = 39
$arr[33,4] = 0
$arr[33,5] = 8
$arr[33,6] = 0
$arr[33,7] = 6
$arr[34,0] = 34
$arr[34,1] = 45391.86111111111
$arr[34,2] = 1
$arr[34,3] = 43
$arr[34,4] = 1
$arr[34,5] = 11
$arr[34,6] = 3
$arr[34,7] = 7
$arr[35,0] = 35
$arr[35,1] = 45391.86805555555
$arr[35,2] = 0
$arr[35,3] = 56
$arr[35,4] = 0
$arr[35,5] = 9
$arr[35,6] = 0
$arr[35,7] = 10

$ws.Range("A2:H37").Value = $arr

$ws.Range("A1").Select()
